# update new orleans xlsx files
#
# 1. Insert a new "State" column into hotel_info (between Hotel_Name and
#    City) and populate it with the hotel's state ("Louisiana").
# 2. Re-order the worksheet tabs so review_info comes before hotel_info.

$wb = $excel.ActiveWorkbook

$hotelSheet = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

# Insert a blank column at C, pushing City/Zip/etc. one column to the right.
$hotelSheet.Columns.Item(3).Insert()

# Fill in the new "State" header + value for the single data row.
$hotelSheet.Range("C1").Value = "State"
$hotelSheet.Range("C2").Value = "Louisiana"

# Move review_info so it becomes the first tab, ahead of hotel_info.
$reviewSheet.Move($hotelSheet)
